$d = $word.ActiveDocument

# delete before[46:52] (1-based)
$startP = $d.Paragraphs.Item(46)
$endP = $d.Paragraphs.Item(52)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()

# replace-resize before[36:44] (1-based) with after[61:65]
$startP = $d.Paragraphs.Item(36)
$endP = $d.Paragraphs.Item(44)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()
$anchor = $d.Paragraphs.Item(35)
$anchor.Range.InsertParagraphAfter()
$insStart = 35 + 1
$d.Paragraphs.Item($insStart).Range.Text = "Impact`r• Discovered systematic race coding errors affecting all Black and Asian-American voters`r• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M`r• Built redistricting platform used by thousands of analysts nationwide`r• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
$d.Paragraphs.Item($insStart + 0).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 1).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 2).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 3).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 4).Range.Style = "Normal"

# replace-resize before[31:34] (1-based) with after[48:59]
$startP = $d.Paragraphs.Item(31)
$endP = $d.Paragraphs.Item(34)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()
$anchor = $d.Paragraphs.Item(30)
$anchor.Range.InsertParagraphAfter()
$insStart = 30 + 1
$d.Paragraphs.Item($insStart).Range.Text = "National Redistricting Platform (2020 - 2021)`rCloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide`rTechnologies: GeoDjango, PostGIS, AWS, Docker, React, Python`rImpact: Reduced mapping costs by 73.5%, saving organizations `$4.7M in operational expenses`rFLEEM Political Polling System (2010 - 2012)`rCompletely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity`rTechnologies: Twilio API, Python, Django, PostgreSQL, JavaScript`rImpact: Saved `$840K in operational costs plus millions in avoided software licensing`rGeospatial Demographic Classification System (2013 - 2016)`rMachine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%`rTechnologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow`rImpact: Corrected demographic data affecting all Black and Asian-American voters nationwide"
$d.Paragraphs.Item($insStart + 0).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 1).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 2).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 3).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 4).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 5).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 6).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 7).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 8).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 9).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 10).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 11).Range.Style = "Normal"

# replace-resize before[15:29] (1-based) with after[44:46]
$startP = $d.Paragraphs.Item(15)
$endP = $d.Paragraphs.Item(29)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()
$anchor = $d.Paragraphs.Item(14)
$anchor.Range.InsertParagraphAfter()
$insStart = 14 + 1
$d.Paragraphs.Item($insStart).Range.Text = "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party`r• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems`r• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"
$d.Paragraphs.Item($insStart + 0).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 1).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 2).Range.Style = "Normal"

# replace-resize before[8:13] (1-based) with after[7:42]
$startP = $d.Paragraphs.Item(8)
$endP = $d.Paragraphs.Item(13)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()
$anchor = $d.Paragraphs.Item(7)
$anchor.Range.InsertParagraphAfter()
$insStart = 7 + 1
$d.Paragraphs.Item($insStart).Range.Text = "Partner - Siege Analytics (Austin, TX) | 2005 - Present`rData, Technology and Strategy Consulting`r• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%`r• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration`r• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%`rData Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023`rCivic Graph & Civic Pulse Director`r• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics`r• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions`r• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture`rAnalytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020`rBig Data Engineering Transformation`r• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS`r• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed`r• Rewrote mission and offerings of department and drafted integration plan with strategy team`rSoftware Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018`rSimCrisis Product Owner/Engineer`r• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement`r• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies`r• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures`rSenior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014`rRACSO Product Owner/Engineer`r• Designed comprehensive survey instruments for specialized voting segments and niche markets`r• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis`r• Wrote RFP and analyzed bids from 1,200 vendors for research platform development`rResearch Director - PCCC (Washington, DC) | 2010 - 2012`rPolitical Research & Data Analysis (FLEEM System)`r• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys`r• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren`r• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver`rSoftware Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011`rGeospatial CRM Development`r• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously`r• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers`r• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill`rProgrammer - Lake Research Partners (Washington, DC) | April 2008 - December 2008"
$d.Paragraphs.Item($insStart + 0).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 1).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 2).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 3).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 4).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 5).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 6).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 7).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 8).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 9).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 10).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 11).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 12).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 13).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 14).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 15).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 16).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 17).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 18).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 19).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 20).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 21).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 22).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 23).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 24).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 25).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 26).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 27).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 28).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 29).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 30).Range.Style = "Heading 3"
$d.Paragraphs.Item($insStart + 31).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 32).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 33).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 34).Range.Style = "Normal"
$d.Paragraphs.Item($insStart + 35).Range.Style = "Heading 3"

# replace before[6:6] (1-based) with after[5:5]
$p = $d.Paragraphs.Item(6)
$rng = $p.Range
[void]$rng.MoveEnd(1, -1)
$rng.Delete()

# replace before[4:4] (1-based) with after[3:3]
$p = $d.Paragraphs.Item(4)
$p.Range.Text = "Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide."

# delete before[2:2] (1-based)
$startP = $d.Paragraphs.Item(2)
$endP = $d.Paragraphs.Item(2)
$r = $d.Range($startP.Range.Start, $endP.Range.End)
$r.Delete()
